$wb = $excel.ActiveWorkbook

# Sheet1 = "BusinessFlow": add new row 5 with TC01_AmazonLogin / login
$ws1 = $wb.Worksheets.Item("BusinessFlow")
$ws1.Range("A5").Value = "TC01_AmazonLogin"
$ws1.Range("B5").Value = "login"
$ws1.Range("A5").Select()

# Sheet2 = "GeneralData": add new row 5 with TC01_AmazonLogin
$ws2 = $wb.Worksheets.Item("GeneralData")
$ws2.Range("A5").Value = "TC01_AmazonLogin"

$ws2.Activate()
$ws2.Range("B5").Select()
